$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.634.32'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.428.19'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +5.25%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '478.72'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +6.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.68'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +13.65%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.995'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.35%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.501'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +7.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.458.59'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +6.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0966'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +11.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.46'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.06%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.323'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +7.06%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.857.31'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '54.804.92'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +4.31%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +8.66%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +13.86%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.454.12'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +5.87%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +9.76%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '313.40'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +5.45%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.76'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +9.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.994'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.41%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +9.48%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '57.00'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.49%  '

$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.165'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +14.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.404'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +10.69%  '

$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.545.06'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +5.32%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +5.63%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0776'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +17.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.996'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '148.43'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.84%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.88'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +6.34%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +9.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.14'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +9.06%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +11.69%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.840'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.57%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '33.04'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.87%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +9.40%  '

$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.598'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +7.03%  '

$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0543'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +7.64%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +9.77%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '254.84'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +26.58%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +11.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0900'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +9.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.930.16'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0221'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +8.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.08'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +9.98%  '
